# Applies the updated crypto price/volume(1h) values and the OKB/VeChain row swap
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $rng = $ws.Range($Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '67.377.97'
Set-TextValue 'E2' '  -2.74%  '
Set-TextValue 'D3' '3.255.62'
Set-TextValue 'E3' '  -5.51%  '
Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  +0.12%  '
Set-TextValue 'D5' '589.62'
Set-TextValue 'E5' '  -3.06%  '
Set-TextValue 'D6' '149.32'
Set-TextValue 'E6' '  -10.67%  '
Set-TextValue 'E7' '  -0.09%  '
Set-TextValue 'D8' '3.251.11'
Set-TextValue 'E8' '  -5.41%  '
Set-TextValue 'D9' '0.542'
Set-TextValue 'E9' '  -8.76%  '
Set-TextValue 'D10' '0.170'
Set-TextValue 'E10' '  -11.15%  '
Set-TextValue 'D11' '6.66'
Set-TextValue 'E11' '  -5.30%  '
Set-TextValue 'D12' '0.503'
Set-TextValue 'E12' '  -10.72%  '
Set-TextValue 'D13' '0.0000245'
Set-TextValue 'E13' '  -8.65%  '
Set-TextValue 'D14' '38.19'
Set-TextValue 'E14' '  -13.81%  '
Set-TextValue 'D15' '3.782.59'
Set-TextValue 'E15' '  -5.64%  '
Set-TextValue 'D16' '67.460.04'
Set-TextValue 'E16' '  -2.81%  '
Set-TextValue 'D17' '3.258.25'
Set-TextValue 'E17' '  -5.61%  '
Set-TextValue 'E18' '  -5.51%  '
Set-TextValue 'D19' '528.18'
Set-TextValue 'E19' '  -8.72%  '
Set-TextValue 'D20' '7.10'
Set-TextValue 'E20' '  -12.99%  '
Set-TextValue 'D21' '14.93'
Set-TextValue 'E21' '  -12.90%  '
Set-TextValue 'D22' '0.753'
Set-TextValue 'E22' '  -10.99%  '
Set-TextValue 'D23' '7.86'
Set-TextValue 'E23' '  -11.61%  '
Set-TextValue 'D24' '85.41'
Set-TextValue 'E24' '  -11.07%  '
Set-TextValue 'D25' '13.45'
Set-TextValue 'E25' '  -11.28%  '
Set-TextValue 'D26' '1.00'
Set-TextValue 'E26' '  +0.03%  '
Set-TextValue 'E27' '  -11.21%  '
Set-TextValue 'D28' '2.14'
Set-TextValue 'E28' '  -11.89%  '
Set-TextValue 'D29' '7.98'
Set-TextValue 'E29' '  -7.49%  '
Set-TextValue 'D30' '29.00'
Set-TextValue 'E30' '  -11.34%  '
Set-TextValue 'E31' '  -3.53%  '
Set-TextValue 'D32' '2.66'
Set-TextValue 'E32' '  -4.95%  '
Set-TextValue 'D33' '6.58'
Set-TextValue 'E33' '  -16.22%  '
Set-TextValue 'D34' '5.69'
Set-TextValue 'E34' '  -13.30%  '
Set-TextValue 'E35' '  -0.13%  '
Set-TextValue 'D36' '512.06'
Set-TextValue 'E36' '  -11.75%  '
Set-TextValue 'B37' 'VeChain'
Set-TextValue 'C37' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D37' '0.0438'
Set-TextValue 'E37' '  -7.01%  '
Set-TextValue 'B38' 'OKB'
Set-TextValue 'C38' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D38' '52.96'
Set-TextValue 'E38' '  -5.54%  '
Set-TextValue 'D39' '0.0851'
Set-TextValue 'D40' '8.91'
Set-TextValue 'E40' '  -15.37%  '
Set-TextValue 'D41' '0.125'
Set-TextValue 'E41' '  -10.61%  '
Set-TextValue 'E42' '  -12.41%  '
Set-TextValue 'D43' '2.925.75'
Set-TextValue 'E43' '  -9.82%  '
Set-TextValue 'D44' '0.265'
Set-TextValue 'E44' '  -10.40%  '
Set-TextValue 'D45' '0.0₃0585'
Set-TextValue 'E45' '  -14.88%  '
Set-TextValue 'D46' '2.18'
Set-TextValue 'E46' '  -9.08%  '
Set-TextValue 'D47' '26.52'
Set-TextValue 'E47' '  -14.64%  '
Set-TextValue 'D49' '2.30'
Set-TextValue 'E49' '  -16.99%  '
Set-TextValue 'E50' '  -10.20%  '
Set-TextValue 'D51' '123.62'
Set-TextValue 'E51' '  -7.78%  '
